# Apply the "Facultades" (departments) list to the "Sitios" sheet, fix the
# selection on both sheets, and widen column B to fit the new long names.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Hoja1": only the active-cell selection changed.
# ---------------------------------------------------------------------
$hoja1 = $wb.Worksheets.Item("Hoja1")
$hoja1.Range("H3").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "Sitios": replace the old single "Comunicación y cultura" row
# with the full list of UCE faculties, clear the now-unused C3/F3
# cells, extend column A with the trailing sequence numbers, widen
# column B, and update the selection.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Sitios")

$ws.Range("B3").Value = "Artes"
$ws.Range("C3").ClearContents()
$ws.Range("F3").ClearContents()

$ws.Range("B4").Value = "Arquitectura y Urbanismo"

$ws.Range("B5").Value = "Ciencias Administrativas"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Ciencias Agrícolas"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Ciencias Económicas"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Ciencias Médicas"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Ciencias Psicológicas"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Ciencias Químicas"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Comunicación Social"

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Cultura Física"

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Ingeniería Ciencias Físicas y Matemática"

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Filosofía, Letras y Ciencias de la Educación"

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "Ingeniería en Geología, Minas, Petróleo y Ambiental"

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Ingeniería Química"

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "Jurisprudencia, Ciencias Políticas y Sociales"

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "Medicina Veterinaria y Zootecnia"

$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "Odontología"

# Row 20 is intentionally left blank (matches source data gap).
$ws.Range("A21").Value = 19
$ws.Range("A22").Value = 20
$ws.Range("A23").Value = 21
$ws.Range("A24").Value = 22
$ws.Range("A25").Value = 23
$ws.Range("A26").Value = 24
$ws.Range("A27").Value = 25

# Widen column B to fit the longer faculty names (stored width "44").
$ws.Range("B1").EntireColumn.ColumnWidth = 43.140625

# Update the saved selection/active cell for this sheet.
$ws.Range("A28").Select() | Out-Null
